$wb = $excel.ActiveWorkbook

# ---- Sheet "gUSD 26.12.24": header row auto-fit height tweak ----
$ws2 = $wb.Worksheets.Item("gUSD 26.12.24")
$ws2.Rows.Item(3).RowHeight = 77.25

# ---- Sheet "mPendle 27.03.25": header row auto-fit height tweak ----
$ws4 = $wb.Worksheets.Item("mPendle 27.03.25")
$ws4.Rows.Item(4).RowHeight = 87.75

# ---- Sheet "gUSD 26.06.25": main data refresh ----
$ws3 = $wb.Worksheets.Item("gUSD 26.06.25")
$ws3.Rows.Item(4).RowHeight = 87.75

# New daily data for rows 23-46 (columns C,D,E,F,G): pos.val, rewards, Impl.APY, undrl.APY 7d, undrl.APY 1d
$cgData = @(
    @(110.225,24.5216,10.22,15.46,10.64),
    @(106.534,25.1289,9.93,14.9,20.2),
    @(106.206,26.8074,9.97,15.52,11.92),
    @(111.667,31.0703,10.6,23.87,73.06),
    @(117.457,33.3026,11.27,26.03,33.86),
    @(120.175,34.4452,11.64,25.47,15.86),
    @(106.816,34.732,10.33,22.51,3.83),
    @(104.64,35.2011,10.18,21.77,6.23),
    @(103.977,35.6981,10.18,20.33,7.58),
    @(103.911,36.4477,10.26,19.45,9.98),
    @(102.681,37.0485,10.21,11.74,8.01),
    @(103.977,38.0318,10.43,9.08,13.37),
    @(103.269,38.4895,10.43,7.72,6.09),
    @(103.248,39.6482,10.51,9.44,16.05),
    @(101.522,40.1837,10.41,9.52,6.74),
    @(99.9045,41.1788,10.31,10.56,14.05),
    @(98.8282,41.4956,10.27,9.71,4.31),
    @(98.9629,42.1785,10.37,9.87,9.1),
    @(107.637,43.4697,11.45,10.49,17.99),
    @(105.329,45.0086,11.27,12.72,22.05),
    @(102.74,45.3913,11.07,11.12,5.02),
    @(100.922,45.9285,10.95,11.16,6.93),
    @(99.8443,46.4583,10.92,10.16,6.94),
    @(99.4208,46.9718,10.97,10.47,6.24)
)
for ($i = 0; $i -lt $cgData.Length; $i++) {
    $r = 23 + $i
    $row = $cgData[$i]
    for ($j = 0; $j -lt 5; $j++) {
        $ws3.Cells.Item($r, 3 + $j).Value = $row[$j]
    }
}

# Extend A column (date, +1 day per row) formula down through row 61
$ws3.Range('A25:A61').Formula = '=A24+1'

# Extend I column (price, C/$A$3) formula down through row 46 (previously only to row 24)
$ws3.Range('I25:I46').Formula = '=C25/$A$3'

# Extend M column (sum of rewards+pos) formula down through row 46 (previously only to row 24)
$ws3.Range('M25:M46').Formula = '=C25+D25'

# ---- Sheet view: select C47 on the active sheet (matches saved cursor position) ----
$ws3.Activate()
$ws3.Range('C47').Select()

